# edit.ps1
# Applies the commit's changes to before.xlsx:
#   1. Sheet "部门情况202404" (dept summary): update F6, G6, H6, I6 values
#      (1100.00 -> 1400.00, 2.00 -> 3.00, 4.59 -> 4.85, 1000.00 -> 1300.00)
#   2. Sheet "经办人情况202404" (agent summary): insert a new data row at
#      row 20 for 金超 / 公司业务四部, pushing the existing rows 20-28 down
#      to rows 21-29 (dimension grows from A1:M28 to A1:M29).
#
# All of the cells involved store their numbers as text (t="inlineStr" in the
# original OOXML), so every write below is done through a leading apostrophe
# (forces text) followed by restoring the cell's original (default) Style so
# that no incidental "@" text-format / style index gets baked into the file.

$wb = $excel.ActiveWorkbook

function Set-TextValue {
    param(
        $Range,
        [string]$Text
    )
    $origStyle = $Range.Style
    $Range.Value = "'" + $Text
    $Range.Style = $origStyle
}

# ---------------------------------------------------------------------
# 1) Sheet "部门情况202404" - row 6 value updates
# ---------------------------------------------------------------------
$wsDept = $wb.Worksheets.Item("部门情况202404")

Set-TextValue $wsDept.Range("F6") "1400.00"
Set-TextValue $wsDept.Range("G6") "3.00"
Set-TextValue $wsDept.Range("H6") "4.85"
Set-TextValue $wsDept.Range("I6") "1300.00"

# ---------------------------------------------------------------------
# 2) Sheet "经办人情况202404" - insert new row 20 (金超 / 公司业务四部)
# ---------------------------------------------------------------------
$wsAgent = $wb.Worksheets.Item("经办人情况202404")

# Shift existing rows 20-28 down to 21-29 and leave row 20 blank.
$wsAgent.Rows(20).Insert()

$newRow = [ordered]@{
    "A20" = "金超"
    "B20" = "公司业务四部"
    "C20" = "150.00"
    "D20" = "1.00"
    "E20" = "0.00"
    "F20" = "0.00"
    "G20" = "150.00"
    "H20" = "1.00"
    "I20" = "0.00"
    "J20" = "0.00"
    "K20" = "0.00"
    "L20" = "400.00"
    "M20" = "2.00"
}

foreach ($addr in $newRow.Keys) {
    Set-TextValue $wsAgent.Range($addr) $newRow[$addr]
}
